$wb = $excel.ActiveWorkbook

# --- Constants sheet: insert 2 new rows above row 23, add the new
#     "line items missing in SAP" exception message, and keep row
#     formatting consistent with the rest of the sheet (Retry Scope for
#     Line items search). ---
$ws = $wb.Worksheets.Item("Constants")
$ws.Rows.Item(23).Resize(2).Insert()

$ws.Range("A23").Value = "ExcpMsg_LineitemsMissing"
$ws.Range("B23").Value = "Couldnot find the line items in SAP"

# Match the row height used throughout the rest of the sheet for the
# newly inserted rows (23 holds the new values, 24/25 stay blank).
$ws.Rows.Item(23).RowHeight = 14.25
$ws.Rows.Item(24).RowHeight = 14.25
$ws.Rows.Item(25).RowHeight = 14.25

# --- Update the remembered selections on each sheet ---
$ws1 = $wb.Worksheets.Item("Settings")
$ws1.Activate()
$ws1.Range("B6").Select()

$ws.Activate()
$ws.Range("B17").Select()
